$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 12.2785232898965
$ws.Range("C2").Value = 4.743621158237222
$ws.Range("D2").Value = 6.753099799763373
$ws.Range("F2").Value = 35.41199441109959
$ws.Range("G2").Value = 3.68519272464414
$ws.Range("K2").Value = 11.71356073839947
$ws.Range("M2").Value = 19.98262175994419
$ws.Range("N2").Value = 21.80415377711843

# Row 3
$ws.Range("B3").Value = 12.06905273458588
$ws.Range("C3").Value = 4.577573048156647
$ws.Range("D3").Value = 6.759434984732452
$ws.Range("F3").Value = 35.19587580391311
$ws.Range("G3").Value = 3.688173457892162
$ws.Range("K3").Value = 11.58343772368868
$ws.Range("M3").Value = 19.39895004617322
$ws.Range("N3").Value = 21.82514085331418

# Row 4
$ws.Range("B4").Value = 11.94241150446019
$ws.Range("C4").Value = 4.474165827493789
$ws.Range("D4").Value = 6.763291768531627
$ws.Range("F4").Value = 35.0702736528912
$ws.Range("G4").Value = 3.690097292112668
$ws.Range("K4").Value = 11.50620355071981
$ws.Range("M4").Value = 19.03987916047495
$ws.Range("N4").Value = 21.83982097483215

# Row 5
$ws.Range("B5").Value = 11.89137578747207
$ws.Range("C5").Value = 4.43173945838295
$ws.Range("D5").Value = 6.764855408377017
$ws.Range("F5").Value = 35.02090867331201
$ws.Range("G5").Value = 3.690904910571271
$ws.Range("K5").Value = 11.47543631172168
$ws.Range("M5").Value = 18.89362315818865
$ws.Range("N5").Value = 21.84625367366913

# Row 6
$ws.Range("B6").Value = 11.88293805331103
$ws.Range("C6").Value = 4.424679661174462
$ws.Range("D6").Value = 6.765114572554237
$ws.Range("F6").Value = 35.01282246827022
$ws.Range("G6").Value = 3.691040445487944
$ws.Range("K6").Value = 11.4703711449807
$ws.Range("M6").Value = 18.86934888505744
$ws.Range("N6").Value = 21.84734900303417

# Row 7
$ws.Range("B7").Value = 11.9417208056379
$ws.Range("C7").Value = 4.47359470505356
$ws.Range("D7").Value = 6.763312888477473
$ws.Range("F7").Value = 35.0696004911551
$ws.Range("G7").Value = 3.6901080880968
$ws.Range("K7").Value = 11.50578570540379
$ws.Range("M7").Value = 19.03790609571439
$ws.Range("N7").Value = 21.83990590557685

# Row 8
$ws.Range("B8").Value = 12.20593319490367
$ws.Range("C8").Value = 4.686714376340799
$ws.Range("D8").Value = 6.755291201188463
$ws.Range("F8").Value = 35.33602573607359
$ws.Range("G8").Value = 3.68620109884472
$ws.Range("K8").Value = 11.66816446911725
$ws.Range("M8").Value = 19.78165806641483
$ws.Range("N8").Value = 21.81101737685132

# Row 9
$ws.Range("B9").Value = 12.73641519088381
$ws.Range("C9").Value = 5.089994578325393
$ws.Range("D9").Value = 6.739284951444473
$ws.Range("F9").Value = 35.91312757518118
$ws.Range("G9").Value = 3.679278437880482
$ws.Range("K9").Value = 12.00601923566829
$ws.Range("M9").Value = 21.2245865646044
$ws.Range("N9").Value = 21.76863243013549

# Row 10
$ws.Range("B10").Value = 13.12931858357977
$ws.Range("C10").Value = 5.373636650518398
$ws.Range("D10").Value = 6.727336878687977
$ws.Range("F10").Value = 36.36803946295402
$ws.Range("G10").Value = 3.674637005933335
$ws.Range("K10").Value = 12.26375242956079
$ws.Range("M10").Value = 22.26274302745233
$ws.Range("N10").Value = 21.74623221569309

# Row 11
$ws.Range("B11").Value = 13.30784894317576
$ws.Range("C11").Value = 5.499268085797598
$ws.Range("D11").Value = 6.721855968312056
$ws.Range("F11").Value = 36.58114411170765
$ws.Range("G11").Value = 3.672620801012331
$ws.Range("K11").Value = 12.38255980803767
$ws.Range("M11").Value = 22.72793157821333
$ws.Range("N11").Value = 21.73794870683227

# Row 12
$ws.Range("B12").Value = 13.37535290111931
$ws.Range("C12").Value = 5.546306011574533
$ws.Range("D12").Value = 6.719773567121205
$ws.Range("F12").Value = 36.66267728726916
$ws.Range("G12").Value = 3.671870911003135
$ws.Range("K12").Value = 12.42773129126066
$ws.Range("M12").Value = 22.90289738513713
$ws.Range("N12").Value = 21.73508679965935

# Row 13
$ws.Range("B13").Value = 13.36082035562351
$ws.Range("C13").Value = 5.536200082812847
$ws.Range("D13").Value = 6.720222361687633
$ws.Range("F13").Value = 36.64508135049545
$ws.Range("G13").Value = 3.67203180962519
$ws.Range("K13").Value = 12.41799540475858
$ws.Range("M13").Value = 22.86527092974791
$ws.Range("N13").Value = 21.73569092584425

# Row 14
$ws.Range("B14").Value = 13.31340497480476
$ws.Range("C14").Value = 5.503148916575229
$ws.Range("D14").Value = 6.721684788119307
$ws.Range("F14").Value = 36.58783545778553
$ws.Range("G14").Value = 3.672558834985058
$ws.Range("K14").Value = 12.3862726969777
$ws.Range("M14").Value = 22.74235092656497
$ws.Range("N14").Value = 21.73770774293635

# Row 15
$ws.Range("B15").Value = 13.28434635479019
$ws.Range("C15").Value = 5.482833017589928
$ws.Range("D15").Value = 6.72257965803778
$ws.Range("F15").Value = 36.55287790058446
$ws.Range("G15").Value = 3.672883422159629
$ws.Range("K15").Value = 12.36686401333856
$ws.Range("M15").Value = 22.66689889642413
$ws.Range("N15").Value = 21.73897892018801

# Row 16
$ws.Range("B16").Value = 13.11764085905617
$ws.Range("C16").Value = 5.365353781829827
$ws.Range("D16").Value = 6.72769412185281
$ws.Range("F16").Value = 36.35423225223779
$ws.Range("G16").Value = 3.674770677972249
$ws.Range("K16").Value = 12.25601585445349
$ws.Range("M16").Value = 22.23218428293447
$ws.Range("N16").Value = 21.7468119882882

# Row 17
$ws.Range("B17").Value = 13.01526794112699
$ws.Range("C17").Value = 5.292378151175222
$ws.Range("D17").Value = 6.730819748485138
$ws.Range("F17").Value = 36.23391220037804
$ws.Range("G17").Value = 3.675952770492777
$ws.Range("K17").Value = 12.18838374305954
$ws.Range("M17").Value = 21.96356023488591
$ws.Range("N17").Value = 21.75210613860436

# Row 18
$ws.Range("B18").Value = 12.95637076143828
$ws.Range("C18").Value = 5.250087040342674
$ws.Range("D18").Value = 6.732613247720066
$ws.Range("F18").Value = 36.16529055162157
$ws.Range("G18").Value = 3.67664164527931
$ws.Range("K18").Value = 12.14963325110551
$ws.Range("M18").Value = 21.8083982812581
$ws.Range("N18").Value = 21.75533059659477

# Row 19
$ws.Range("B19").Value = 12.93642886752519
$ws.Range("C19").Value = 5.235715005545777
$ws.Range("D19").Value = 6.733219770861975
$ws.Range("F19").Value = 36.14215814748193
$ws.Range("G19").Value = 3.6768764291045
$ws.Range("K19").Value = 12.13654003073058
$ws.Range("M19").Value = 21.75575591982339
$ws.Range("N19").Value = 21.75645313259932

# Row 20
$ws.Range("B20").Value = 13.02616779751237
$ws.Range("C20").Value = 5.30017974745015
$ws.Range("D20").Value = 6.730487465385283
$ws.Range("F20").Value = 36.24666047517462
$ws.Range("G20").Value = 3.675826007368153
$ws.Range("K20").Value = 12.19556811552417
$ws.Range("M20").Value = 21.99222498045038
$ws.Range("N20").Value = 21.75152399441235

# Row 21
$ws.Range("B21").Value = 13.3273353412027
$ws.Range("C21").Value = 5.512871740055999
$ws.Range("D21").Value = 6.72125542795661
$ws.Range("F21").Value = 36.60462771616771
$ws.Range("G21").Value = 3.67240366637474
$ws.Range("K21").Value = 12.39558583533031
$ws.Range("M21").Value = 22.77848913762698
$ws.Range("N21").Value = 21.7371078890271

# Row 22
$ws.Range("B22").Value = 13.52354206943044
$ws.Range("C22").Value = 5.648736156536179
$ws.Range("D22").Value = 6.715181382973158
$ws.Range("F22").Value = 36.8434226350095
$ws.Range("G22").Value = 3.670246217153091
$ws.Range("K22").Value = 12.52734843347753
$ws.Range("M22").Value = 23.28534712016935
$ws.Range("N22").Value = 21.72928865427492

# Row 23
$ws.Range("B23").Value = 13.41890303084811
$ws.Range("C23").Value = 5.576524437619439
$ws.Range("D23").Value = 6.718427021162299
$ws.Range("F23").Value = 36.71554748997669
$ws.Range("G23").Value = 3.671390466010058
$ws.Range("K23").Value = 12.45694327347124
$ws.Range("M23").Value = 23.01552215159554
$ws.Range("N23").Value = 21.73331505955697

# Row 24
$ws.Range("B24").Value = 13.02124009657133
$ws.Range("C24").Value = 5.296653690093194
$ws.Range("D24").Value = 6.73063770134319
$ws.Range("F24").Value = 36.24089525707272
$ws.Range("G24").Value = 3.675883288083268
$ws.Range("K24").Value = 12.19231964596887
$ws.Range("M24").Value = 21.97926789822037
$ws.Range("N24").Value = 21.75178661876633

# Row 25
$ws.Range("B25").Value = 12.59203940423148
$ws.Range("C25").Value = 4.982873650075696
$ws.Range("D25").Value = 6.743646645815053
$ws.Range("F25").Value = 35.75139385895853
$ws.Range("G25").Value = 3.681072696312212
$ws.Range("K25").Value = 11.9127786874237
$ws.Range("M25").Value = 20.83724926061207
$ws.Range("N25").Value = 21.77856691566708
